$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2-21 for columns I (I0) and J (IF)
$data = @(
    @(1, 5),
    @(1, 5),
    @(1, 3),
    @(1, 5),
    @(1, 6),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 5),
    @(1, 4),
    @(1, 5),
    @(2, 6),
    @(5, 8),
    @(5, 8),
    @(9, 9),
    @(1, 3),
    @(6, 6)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $data[$idx][0]
    $ws.Cells.Item($row, 10).Value = $data[$idx][1]
}
